# Generate Report for Handoff
#
# Adds two new tracked files (70fbc62f-...md and d0826aa6-...md) to the
# localization-status report. Each of the three sheets (Overview, zh-cn,
# de-de) gets two new rows inserted right before the ".localization-config"
# row, which is pushed down to make room.

$wb = $excel.ActiveWorkbook

# Re-usable constants
$colorLink = 15570276   # RGB(100,149,237) == FF6495ED, matches the workbook's HyperLink style
$underlineSingle = 2    # xlUnderlineStyleSingle

function Style-AsLink($rng) {
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
    $rng.Font.Underline = $underlineSingle
    $rng.Font.Color = $colorLink
}

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Rows.Item(4).Insert()
$ws1.Rows.Item(4).Insert()

$ws1.Range("A4").Value2 = "70fbc62f-fe8d-4b37-94ee-52699469a693.md"
$ws1.Range("B4").Value2 = "Ready for handoff"
$ws1.Range("C4").Value2 = "Ready for handoff"

$ws1.Range("A5").Value2 = "d0826aa6-49cb-413f-b3b3-1728b42d77b2.md"
$ws1.Range("B5").Value2 = "Ready for handoff"
$ws1.Range("C5").Value2 = "Ready for handoff"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/888c6be3f98264881231243a34f4d7647df82067/e2e/022138e0-8300-4b90-a5cc-691146087493.md", [Type]::Missing, [Type]::Missing, "022138e0-8300-4b90-a5cc-691146087493.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/07324d3385ddf113621b306083514f9ccfcf607d/e2e/0460c0ac-d1f2-4c14-b2fe-992dc46fcb23.md", [Type]::Missing, [Type]::Missing, "0460c0ac-d1f2-4c14-b2fe-992dc46fcb23.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/888c6be3f98264881231243a34f4d7647df82067/e2e/70fbc62f-fe8d-4b37-94ee-52699469a693.md", [Type]::Missing, [Type]::Missing, "70fbc62f-fe8d-4b37-94ee-52699469a693.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/888c6be3f98264881231243a34f4d7647df82067/e2e/d0826aa6-49cb-413f-b3b3-1728b42d77b2.md", [Type]::Missing, [Type]::Missing, "d0826aa6-49cb-413f-b3b3-1728b42d77b2.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/888c6be3f98264881231243a34f4d7647df82067/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

Style-AsLink($ws1.Range("A2"))
Style-AsLink($ws1.Range("A3"))
Style-AsLink($ws1.Range("A4"))
Style-AsLink($ws1.Range("A5"))
Style-AsLink($ws1.Range("A6"))

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows.Item(4).Insert()
$ws2.Rows.Item(4).Insert()

$ws2.Range("A4").Value2 = "70fbc62f-fe8d-4b37-94ee-52699469a693.md"
$ws2.Range("B4").Value2 = "Ready for handoff"
$ws2.Range("C4").Value2 = "70fbc62f-fe8d-4b37-94ee-52699469a693.34ad18ae0cd84a22a7726407fc2b5359b366fc60.zh-cn.xlf"
$ws2.Range("D4").Value2 = "2016-03-09 04:47:20"
$ws2.Range("G4").Value2 = "0001-01-01 00:00:00"
$ws2.Range("H4").Value2 = "Include"

$ws2.Range("A5").Value2 = "d0826aa6-49cb-413f-b3b3-1728b42d77b2.md"
$ws2.Range("B5").Value2 = "Ready for handoff"
$ws2.Range("C5").Value2 = "d0826aa6-49cb-413f-b3b3-1728b42d77b2.8888820a040029c4899bb3a7870def3e041b618d.zh-cn.xlf"
$ws2.Range("D5").Value2 = "2016-03-09 04:47:20"
$ws2.Range("G5").Value2 = "0001-01-01 00:00:00"
$ws2.Range("H5").Value2 = "Include"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/888c6be3f98264881231243a34f4d7647df82067/e2e/022138e0-8300-4b90-a5cc-691146087493.md", [Type]::Missing, [Type]::Missing, "022138e0-8300-4b90-a5cc-691146087493.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e3264f1700de2634665ccf978071e2c7e851953c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/022138e0-8300-4b90-a5cc-691146087493.a2a2cfdf0ed001351966bdfc3eb713672368944d.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "022138e0-8300-4b90-a5cc-691146087493.a2a2cfdf0ed001351966bdfc3eb713672368944d.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/056d7de94347f23e15c0ab0566e69fc5d4e4ee21/e2e/022138e0-8300-4b90-a5cc-691146087493.md", [Type]::Missing, [Type]::Missing, "022138e0-8300-4b90-a5cc-691146087493.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dff0b119bebbef85eb6e710493bc8ea5a044c3de/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/022138e0-8300-4b90-a5cc-691146087493.a2a2cfdf0ed001351966bdfc3eb713672368944d.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "022138e0-8300-4b90-a5cc-691146087493.a2a2cfdf0ed001351966bdfc3eb713672368944d.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/07324d3385ddf113621b306083514f9ccfcf607d/e2e/0460c0ac-d1f2-4c14-b2fe-992dc46fcb23.md", [Type]::Missing, [Type]::Missing, "0460c0ac-d1f2-4c14-b2fe-992dc46fcb23.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4339236898bcb017f41c8d7a27d892bdfd0407ba/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0460c0ac-d1f2-4c14-b2fe-992dc46fcb23.5d9ec070339fe1467f90a34cf9f32c4dbcf1a1ae.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "0460c0ac-d1f2-4c14-b2fe-992dc46fcb23.5d9ec070339fe1467f90a34cf9f32c4dbcf1a1ae.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/888c6be3f98264881231243a34f4d7647df82067/e2e/70fbc62f-fe8d-4b37-94ee-52699469a693.md", [Type]::Missing, [Type]::Missing, "70fbc62f-fe8d-4b37-94ee-52699469a693.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/34ad18ae0cd84a22a7726407fc2b5359b366fc60/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/70fbc62f-fe8d-4b37-94ee-52699469a693.34ad18ae0cd84a22a7726407fc2b5359b366fc60.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "70fbc62f-fe8d-4b37-94ee-52699469a693.34ad18ae0cd84a22a7726407fc2b5359b366fc60.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/888c6be3f98264881231243a34f4d7647df82067/e2e/d0826aa6-49cb-413f-b3b3-1728b42d77b2.md", [Type]::Missing, [Type]::Missing, "d0826aa6-49cb-413f-b3b3-1728b42d77b2.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8888820a040029c4899bb3a7870def3e041b618d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d0826aa6-49cb-413f-b3b3-1728b42d77b2.8888820a040029c4899bb3a7870def3e041b618d.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "d0826aa6-49cb-413f-b3b3-1728b42d77b2.8888820a040029c4899bb3a7870def3e041b618d.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/888c6be3f98264881231243a34f4d7647df82067/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

Style-AsLink($ws2.Range("A2"))
Style-AsLink($ws2.Range("C2"))
Style-AsLink($ws2.Range("E2"))
Style-AsLink($ws2.Range("F2"))
Style-AsLink($ws2.Range("A3"))
Style-AsLink($ws2.Range("C3"))
Style-AsLink($ws2.Range("A4"))
Style-AsLink($ws2.Range("C4"))
Style-AsLink($ws2.Range("A5"))
Style-AsLink($ws2.Range("C5"))
Style-AsLink($ws2.Range("A6"))

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows.Item(4).Insert()
$ws3.Rows.Item(4).Insert()

$ws3.Range("A4").Value2 = "70fbc62f-fe8d-4b37-94ee-52699469a693.md"
$ws3.Range("B4").Value2 = "Ready for handoff"
$ws3.Range("C4").Value2 = "70fbc62f-fe8d-4b37-94ee-52699469a693.34ad18ae0cd84a22a7726407fc2b5359b366fc60.de-de.xlf"
$ws3.Range("D4").Value2 = "2016-03-09 04:47:22"
$ws3.Range("G4").Value2 = "0001-01-01 00:00:00"
$ws3.Range("H4").Value2 = "Include"

$ws3.Range("A5").Value2 = "d0826aa6-49cb-413f-b3b3-1728b42d77b2.md"
$ws3.Range("B5").Value2 = "Ready for handoff"
$ws3.Range("C5").Value2 = "d0826aa6-49cb-413f-b3b3-1728b42d77b2.8888820a040029c4899bb3a7870def3e041b618d.de-de.xlf"
$ws3.Range("D5").Value2 = "2016-03-09 04:47:22"
$ws3.Range("G5").Value2 = "0001-01-01 00:00:00"
$ws3.Range("H5").Value2 = "Include"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/888c6be3f98264881231243a34f4d7647df82067/e2e/022138e0-8300-4b90-a5cc-691146087493.md", [Type]::Missing, [Type]::Missing, "022138e0-8300-4b90-a5cc-691146087493.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c53ae1f66925475b96065a9d6c25e50a3322b7c4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/022138e0-8300-4b90-a5cc-691146087493.a2a2cfdf0ed001351966bdfc3eb713672368944d.de-de.xlf", [Type]::Missing, [Type]::Missing, "022138e0-8300-4b90-a5cc-691146087493.a2a2cfdf0ed001351966bdfc3eb713672368944d.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/208367b039c81f896e0ac65827902d7e3c8c3474/e2e/022138e0-8300-4b90-a5cc-691146087493.md", [Type]::Missing, [Type]::Missing, "022138e0-8300-4b90-a5cc-691146087493.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8c63f14b043400fb96f2d8a1e44294d6aa10f86e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/022138e0-8300-4b90-a5cc-691146087493.a2a2cfdf0ed001351966bdfc3eb713672368944d.de-de.xlf", [Type]::Missing, [Type]::Missing, "022138e0-8300-4b90-a5cc-691146087493.a2a2cfdf0ed001351966bdfc3eb713672368944d.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/07324d3385ddf113621b306083514f9ccfcf607d/e2e/0460c0ac-d1f2-4c14-b2fe-992dc46fcb23.md", [Type]::Missing, [Type]::Missing, "0460c0ac-d1f2-4c14-b2fe-992dc46fcb23.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6d055b294614f4cd9f77e35d78ffbf8db70a92c1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0460c0ac-d1f2-4c14-b2fe-992dc46fcb23.5d9ec070339fe1467f90a34cf9f32c4dbcf1a1ae.de-de.xlf", [Type]::Missing, [Type]::Missing, "0460c0ac-d1f2-4c14-b2fe-992dc46fcb23.5d9ec070339fe1467f90a34cf9f32c4dbcf1a1ae.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/888c6be3f98264881231243a34f4d7647df82067/e2e/70fbc62f-fe8d-4b37-94ee-52699469a693.md", [Type]::Missing, [Type]::Missing, "70fbc62f-fe8d-4b37-94ee-52699469a693.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/34ad18ae0cd84a22a7726407fc2b5359b366fc60/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/70fbc62f-fe8d-4b37-94ee-52699469a693.34ad18ae0cd84a22a7726407fc2b5359b366fc60.de-de.xlf", [Type]::Missing, [Type]::Missing, "70fbc62f-fe8d-4b37-94ee-52699469a693.34ad18ae0cd84a22a7726407fc2b5359b366fc60.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/888c6be3f98264881231243a34f4d7647df82067/e2e/d0826aa6-49cb-413f-b3b3-1728b42d77b2.md", [Type]::Missing, [Type]::Missing, "d0826aa6-49cb-413f-b3b3-1728b42d77b2.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8888820a040029c4899bb3a7870def3e041b618d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d0826aa6-49cb-413f-b3b3-1728b42d77b2.8888820a040029c4899bb3a7870def3e041b618d.de-de.xlf", [Type]::Missing, [Type]::Missing, "d0826aa6-49cb-413f-b3b3-1728b42d77b2.8888820a040029c4899bb3a7870def3e041b618d.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/888c6be3f98264881231243a34f4d7647df82067/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

Style-AsLink($ws3.Range("A2"))
Style-AsLink($ws3.Range("C2"))
Style-AsLink($ws3.Range("E2"))
Style-AsLink($ws3.Range("F2"))
Style-AsLink($ws3.Range("A3"))
Style-AsLink($ws3.Range("C3"))
Style-AsLink($ws3.Range("A4"))
Style-AsLink($ws3.Range("C4"))
Style-AsLink($ws3.Range("A5"))
Style-AsLink($ws3.Range("C5"))
Style-AsLink($ws3.Range("A6"))
